$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows referencing the removed "Inflammatory-Mac" target cluster
# (delete from the bottom up so row indices of earlier rows stay stable)
$ws.Rows(9).Delete()
$ws.Rows(4).Delete()

# Refresh remaining rows with updated TPM-derived values
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Wnt1"
$ws.Cells.Item(2, 3).Value = "Fzd2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.027123
$ws.Cells.Item(2, 8).Value = 0.081369
$ws.Cells.Item(2, 9).Value = 0.0960827240265261
$ws.Cells.Item(2, 10).Value = 0.09608272402652611
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.08425766666666668
$ws.Cells.Item(2, 14).Value = 0.252773
$ws.Cells.Item(2, 15).Value = 0.007654801123801229
$ws.Cells.Item(2, 16).Value = 0.008027752567511702
$ws.Cells.Item(2, 17).Value = 0.002285320693
$ws.Cells.Item(2, 18).Value = 0.020567886237
$ws.Cells.Item(2, 19).Value = 0.0007354941438561353
$ws.Cells.Item(2, 20).Value = 0.0007713283344974632
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Wnt1"
$ws.Cells.Item(3, 3).Value = "Fzd2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.027123
$ws.Cells.Item(3, 8).Value = 0.081369
$ws.Cells.Item(3, 9).Value = 0.0960827240265261
$ws.Cells.Item(3, 10).Value = 0.09608272402652611
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 9.380691666666667
$ws.Cells.Item(3, 14).Value = 28.142075
$ws.Cells.Item(3, 15).Value = 0.8522349591772004
$ws.Cells.Item(3, 16).Value = 0.8937569077249424
$ws.Cells.Item(3, 17).Value = 0.254432500075
$ws.Cells.Item(3, 18).Value = 2.289892500675
$ws.Cells.Item(3, 19).Value = 0.08188505638838067
$ws.Cells.Item(3, 20).Value = 0.08587459831173701
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Wnt1"
$ws.Cells.Item(4, 3).Value = "Fzd2"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.027123
$ws.Cells.Item(4, 8).Value = 0.081369
$ws.Cells.Item(4, 9).Value = 0.0960827240265261
$ws.Cells.Item(4, 10).Value = 0.09608272402652611
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.534105
$ws.Cells.Item(4, 14).Value = 3.06821
$ws.Cells.Item(4, 15).Value = 0.139373295542195
$ws.Cells.Item(4, 16).Value = 0.09744249071366434
$ws.Cells.Item(4, 17).Value = 0.041609529915
$ws.Cells.Item(4, 18).Value = 0.24965717949
$ws.Cells.Item(4, 19).Value = 0.01339136589224818
$ws.Cells.Item(4, 20).Value = 0.009362539943698345
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Wnt1"
$ws.Cells.Item(5, 3).Value = "Fzd2"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.027123
$ws.Cells.Item(5, 8).Value = 0.081369
$ws.Cells.Item(5, 9).Value = 0.0960827240265261
$ws.Cells.Item(5, 10).Value = 0.09608272402652611
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.008111666666666666
$ws.Cells.Item(5, 14).Value = 0.024335
$ws.Cells.Item(5, 15).Value = 0.0007369441568035466
$ws.Cells.Item(5, 16).Value = 0.0007728489938814559
$ws.Cells.Item(5, 17).Value = 0.000220012735
$ws.Cells.Item(5, 18).Value = 0.001980114615
$ws.Cells.Item(5, 19).Value = [double]"7.080760204111614E-05"
$ws.Cells.Item(5, 20).Value = [double]"7.425743659329029E-05"
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Wnt1"
$ws.Cells.Item(6, 3).Value = "Fzd2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.255165
$ws.Cells.Item(6, 8).Value = 0.7654949999999999
$ws.Cells.Item(6, 9).Value = 0.9039172759734738
$ws.Cells.Item(6, 10).Value = 0.9039172759734738
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.08425766666666668
$ws.Cells.Item(6, 14).Value = 0.252773
$ws.Cells.Item(6, 15).Value = 0.007654801123801229
$ws.Cells.Item(6, 16).Value = 0.008027752567511702
$ws.Cells.Item(6, 17).Value = 0.021499607515
$ws.Cells.Item(6, 18).Value = 0.193496467635
$ws.Cells.Item(6, 19).Value = 0.006919306979945093
$ws.Cells.Item(6, 20).Value = 0.007256424233014238
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Wnt1"
$ws.Cells.Item(7, 3).Value = "Fzd2"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.255165
$ws.Cells.Item(7, 8).Value = 0.7654949999999999
$ws.Cells.Item(7, 9).Value = 0.9039172759734738
$ws.Cells.Item(7, 10).Value = 0.9039172759734738
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 9.380691666666667
$ws.Cells.Item(7, 14).Value = 28.142075
$ws.Cells.Item(7, 15).Value = 0.8522349591772004
$ws.Cells.Item(7, 16).Value = 0.8937569077249424
$ws.Cells.Item(7, 17).Value = 2.393624189125
$ws.Cells.Item(7, 18).Value = 21.542617702125
$ws.Cells.Item(7, 19).Value = 0.7703499027888197
$ws.Cells.Item(7, 20).Value = 0.8078823094132054
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Wnt1"
$ws.Cells.Item(8, 3).Value = "Fzd2"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.255165
$ws.Cells.Item(8, 8).Value = 0.7654949999999999
$ws.Cells.Item(8, 9).Value = 0.9039172759734738
$ws.Cells.Item(8, 10).Value = 0.9039172759734738
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.534105
$ws.Cells.Item(8, 14).Value = 3.06821
$ws.Cells.Item(8, 15).Value = 0.139373295542195
$ws.Cells.Item(8, 16).Value = 0.09744249071366434
$ws.Cells.Item(8, 17).Value = 0.391449902325
$ws.Cells.Item(8, 18).Value = 2.34869941395
$ws.Cells.Item(8, 19).Value = 0.1259819296499468
$ws.Cells.Item(8, 20).Value = 0.088079950769966
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Wnt1"
$ws.Cells.Item(9, 3).Value = "Fzd2"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.255165
$ws.Cells.Item(9, 8).Value = 0.7654949999999999
$ws.Cells.Item(9, 9).Value = 0.9039172759734738
$ws.Cells.Item(9, 10).Value = 0.9039172759734738
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.008111666666666666
$ws.Cells.Item(9, 14).Value = 0.024335
$ws.Cells.Item(9, 15).Value = 0.0007369441568035466
$ws.Cells.Item(9, 16).Value = 0.0007728489938814559
$ws.Cells.Item(9, 17).Value = 0.002069813425
$ws.Cells.Item(9, 18).Value = 0.018628320825
$ws.Cells.Item(9, 19).Value = 0.0006661365547624304
$ws.Cells.Item(9, 20).Value = 0.0006985915572881655
